# Lecture partielle de l'EDT M1 MIAGE.
# Shift every class-week date forward by 1096 days (3 years later, same
# month/day) and refresh the French weekday label next to each date so it
# matches the new date. Also fix a stray time typo (16:45 -> 15:45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Index 0 = lundi ... 6 = dimanche. For an Excel date serial number
# (days since 1899-12-30), dayNames[(serial + 5) % 7] gives the French
# weekday name.
$dayNames = @("lundi", "mardi", "mercredi", "jeudi", "vendredi", "samedi", "dimanche")

$dayShift = 1096
$dateRows = @(2, 5, 8, 11, 14, 16, 19, 22, 25, 28, 31, 34)

foreach ($r in $dateRows) {
    $oldSerial = $ws.Cells.Item($r, 1).Value2
    $newSerial = $oldSerial + $dayShift
    $ws.Cells.Item($r, 1).Value = $newSerial

    $idx = ($newSerial + 5) % 7
    $dayName = $dayNames[$idx]
    $ws.Cells.Item($r, 2).Value = $dayName
}

# Correction de l'heure : 16:45 -> 15:45
$ws.Cells.Item(6, 4).Value = "15:45"
